$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) from the last existing data row (A4) onto the
# new label cells in column A before filling in their values.
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("A5").Value = "spectral_1"
$ws.Range("B5").Value = -0.3862637995641255
$ws.Range("C5").Value = 0.4725908679043076
$ws.Range("D5").Value = 4.155297896218482

$ws.Range("A6").Value = "hierarchical_1"
$ws.Range("B6").Value = 0.9939182672879472
$ws.Range("C6").Value = 36211.61359472208
$ws.Range("D6").Value = 0.2206786182786379
